$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 (Hydrogen, Iron & steel) value
$ws.Range("B3").Value = 1493901.880850142

# Clear D3 value (Hydrogen, Non-metallic minerals) -> becomes empty cell
$ws.Range("D3").ClearContents()

# Update C4 (Methanol, Chemicals) value
$ws.Range("C4").Value = 4314.65834554362

# Update C5 (Ammonia, Chemicals) value
$ws.Range("C5").Value = 11414.72293445124

# Row 7: rename "Other" -> "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 310.5333032218452

# Row 8 (new): "Other" row with value in D8, copying formatting from row 7
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 24.50151034028678
